# Update Reln-Vldlr LR-pair sheet with refreshed TPM-based NATMI values.
# The "Inflammatory-Mac" target-cluster category is dropped, collapsing
# each sending cluster's 4 remaining target rows (ECs/FAPs/MuSCs/Resolving-Mac)
# into rows 2-13, and the trailing now-unused rows 14-16 are removed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Reln"
$ws.Range("C2").Value = "Vldlr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.05089966666666667
$ws.Range("H2").Value = 0.152699
$ws.Range("I2").Value = 0.02671091810242436
$ws.Range("J2").Value = 0.03728162213961778
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.104012
$ws.Range("N2").Value = 0.312036
$ws.Range("O2").Value = 0.01457989208339885
$ws.Range("P2").Value = 0.01622836922145579
$ws.Range("Q2").Value = 0.005294176129333333
$ws.Range("R2").Value = 0.047647585164
$ws.Range("S2").Value = 0.0003894423033818521
$ws.Range("T2").Value = 0.0006050199292565181

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Reln"
$ws.Range("C3").Value = "Vldlr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.05089966666666667
$ws.Range("H3").Value = 0.152699
$ws.Range("I3").Value = 0.02671091810242436
$ws.Range("J3").Value = 0.03728162213961778
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.808482333333333
$ws.Range("N3").Value = 14.425447
$ws.Range("O3").Value = 0.674029472608256
$ws.Range("P3").Value = 0.7502386907297295
$ws.Range("Q3").Value = 0.2447501479392222
$ws.Range("R3").Value = 2.202751331453
$ws.Range("S3").Value = 0.01800394604145941
$ws.Range("T3").Value = 0.02797011538230734

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Reln"
$ws.Range("C4").Value = "Vldlr"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.05089966666666667
$ws.Range("H4").Value = 0.152699
$ws.Range("I4").Value = 0.02671091810242436
$ws.Range("J4").Value = 0.03728162213961778
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.1739945
$ws.Range("N4").Value = 4.347989
$ws.Range("O4").Value = 0.3047398877043289
$ws.Range("P4").Value = 0.2261302249190105
$ws.Range("Q4").Value = 0.1106555953851667
$ws.Range("R4").Value = 0.663933572311
$ws.Range("S4").Value = 0.008139882183012325
$ws.Range("T4").Value = 0.008430501599777332

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Reln"
$ws.Range("C5").Value = "Vldlr"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.05089966666666667
$ws.Range("H5").Value = 0.152699
$ws.Range("I5").Value = 0.02671091810242436
$ws.Range("J5").Value = 0.03728162213961778
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.047446
$ws.Range("N5").Value = 0.142338
$ws.Range("O5").Value = 0.006650747604016287
$ws.Range("P5").Value = 0.007402715129804173
$ws.Range("Q5").Value = 0.002414985584666666
$ws.Range("R5").Value = 0.021734870262
$ws.Range("S5").Value = 0.0001776475745707741
$ws.Range("T5").Value = 0.0002759852282765908

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Reln"
$ws.Range("C6").Value = "Vldlr"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.2337746666666667
$ws.Range("H6").Value = 0.7013240000000001
$ws.Range("I6").Value = 0.1226793098007496
$ws.Range("J6").Value = 0.1712289953794413
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.104012
$ws.Range("N6").Value = 0.312036
$ws.Range("O6").Value = 0.01457989208339885
$ws.Range("P6").Value = 0.01622836922145579
$ws.Range("Q6").Value = 0.02431537062933333
$ws.Range("R6").Value = 0.218838335664
$ws.Range("S6").Value = 0.001788651097760785
$ws.Range("T6").Value = 0.002778767358436521

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Reln"
$ws.Range("C7").Value = "Vldlr"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.2337746666666667
$ws.Range("H7").Value = 0.7013240000000001
$ws.Range("I7").Value = 0.1226793098007496
$ws.Range("J7").Value = 0.1712289953794413
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.808482333333333
$ws.Range("N7").Value = 14.425447
$ws.Range("O7").Value = 0.674029472608256
$ws.Range("P7").Value = 0.7502386907297295
$ws.Range("Q7").Value = 1.124101354647556
$ws.Range("R7").Value = 10.116912191828
$ws.Range("S7").Value = 0.0826894704849441
$ws.Range("T7").Value = 0.1284626173084389

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Reln"
$ws.Range("C8").Value = "Vldlr"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.2337746666666667
$ws.Range("H8").Value = 0.7013240000000001
$ws.Range("I8").Value = 0.1226793098007496
$ws.Range("J8").Value = 0.1712289953794413
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.1739945
$ws.Range("N8").Value = 4.347989
$ws.Range("O8").Value = 0.3047398877043289
$ws.Range("P8").Value = 0.2261302249190105
$ws.Range("Q8").Value = 0.5082248395726667
$ws.Range("R8").Value = 3.049349037436
$ws.Range("S8").Value = 0.03738527909232501
$ws.Range("T8").Value = 0.03872005123780927

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Reln"
$ws.Range("C9").Value = "Vldlr"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.2337746666666667
$ws.Range("H9").Value = 0.7013240000000001
$ws.Range("I9").Value = 0.1226793098007496
$ws.Range("J9").Value = 0.1712289953794413
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.047446
$ws.Range("N9").Value = 0.142338
$ws.Range("O9").Value = 0.006650747604016287
$ws.Range("P9").Value = 0.007402715129804173
$ws.Range("Q9").Value = 0.01109167283466667
$ws.Range("R9").Value = 0.09982505551200001
$ws.Range("S9").Value = 0.0008159091257197072
$ws.Range("T9").Value = 0.001267559474756559

$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Reln"
$ws.Range("C10").Value = "Vldlr"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.620901
$ws.Range("H10").Value = 3.241802
$ws.Range("I10").Value = 0.8506097720968261
$ws.Range("J10").Value = 0.791489382480941
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.104012
$ws.Range("N10").Value = 0.312036
$ws.Range("O10").Value = 0.01457989208339885
$ws.Range("P10").Value = 0.01622836922145579
$ws.Range("Q10").Value = 0.168593154812
$ws.Range("R10").Value = 1.011558928872
$ws.Range("S10").Value = 0.01240179868225622
$ws.Range("T10").Value = 0.01284458193376276

$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Reln"
$ws.Range("C11").Value = "Vldlr"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.620901
$ws.Range("H11").Value = 3.241802
$ws.Range("I11").Value = 0.8506097720968261
$ws.Range("J11").Value = 0.791489382480941
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 4.808482333333333
$ws.Range("N11").Value = 14.425447
$ws.Range("O11").Value = 0.674029472608256
$ws.Range("P11").Value = 0.7502386907297295
$ws.Range("Q11").Value = 7.794073822582333
$ws.Range("R11").Value = 46.764442935494
$ws.Range("S11").Value = 0.5733360560818525
$ws.Range("T11").Value = 0.5938059580389833

$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Reln"
$ws.Range("C12").Value = "Vldlr"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.620901
$ws.Range("H12").Value = 3.241802
$ws.Range("I12").Value = 0.8506097720968261
$ws.Range("J12").Value = 0.791489382480941
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 2.1739945
$ws.Range("N12").Value = 4.347989
$ws.Range("O12").Value = 0.3047398877043289
$ws.Range("P12").Value = 0.2261302249190105
$ws.Range("Q12").Value = 3.5238298590445
$ws.Range("R12").Value = 14.095319436178
$ws.Range("S12").Value = 0.2592147264289916
$ws.Range("T12").Value = 0.178979672081424

$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Reln"
$ws.Range("C13").Value = "Vldlr"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.620901
$ws.Range("H13").Value = 3.241802
$ws.Range("I13").Value = 0.8506097720968261
$ws.Range("J13").Value = 0.791489382480941
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.047446
$ws.Range("N13").Value = 0.142338
$ws.Range("O13").Value = 0.006650747604016287
$ws.Range("P13").Value = 0.007402715129804173
$ws.Range("Q13").Value = 0.07690526884599999
$ws.Range("R13").Value = 0.4614316130759999
$ws.Range("S13").Value = 0.005657190903725806
$ws.Range("T13").Value = 0.005859170426771024

$ws.Range("A14:T16").EntireRow.Delete()
